# Weekly update: prepend the newest week's "Apio" (Vega Central Mapocho de
# Santiago) price records. Two new data rows are inserted above the current
# first data block (row 529), pushing all the existing rows down by two.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for the two new rows; inserting at 529 twice shifts old 529.. down to 531..
$ws.Rows.Item(529).Insert()
$ws.Rows.Item(529).Insert()

# Copy the (now shifted) following row's formatting onto the two fresh rows so
# the date column keeps its existing custom number format, etc.
$ws.Range("A531:R532").Copy()
$ws.Range("A529:R530").PasteSpecial(-4122)

# --- New row 529: Apio, Americana (o), Primera ---
$ws.Cells.Item(529, 1).Value = 9
$ws.Cells.Item(529, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(529, 3).Value = "Metropolitana"
$ws.Cells.Item(529, 4).Value = "2023-10-24"
$ws.Cells.Item(529, 5).Value = 13
$ws.Cells.Item(529, 6).Value = 100112017
$ws.Cells.Item(529, 7).Value = "Apio"
$ws.Cells.Item(529, 8).Value = "Americana (o)"
$ws.Cells.Item(529, 9).Value = "Primera"
$ws.Cells.Item(529, 10).Value = 70
$ws.Cells.Item(529, 11).Value = 6000
$ws.Cells.Item(529, 12).Value = 7000
$ws.Cells.Item(529, 13).Value = 6514
$ws.Cells.Item(529, 14).Value = "`$/docena de matas"
$ws.Cells.Item(529, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(529, 16).Value = 1086
$ws.Cells.Item(529, 17).Value = 6
$ws.Cells.Item(529, 18).Value = "Hortaliza"

# --- New row 530: Apio, Americana (o), Segunda ---
$ws.Cells.Item(530, 1).Value = 9
$ws.Cells.Item(530, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(530, 3).Value = "Metropolitana"
$ws.Cells.Item(530, 4).Value = "2023-10-24"
$ws.Cells.Item(530, 5).Value = 13
$ws.Cells.Item(530, 6).Value = 100112017
$ws.Cells.Item(530, 7).Value = "Apio"
$ws.Cells.Item(530, 8).Value = "Americana (o)"
$ws.Cells.Item(530, 9).Value = "Segunda"
$ws.Cells.Item(530, 10).Value = 52
$ws.Cells.Item(530, 11).Value = 5000
$ws.Cells.Item(530, 12).Value = 5000
$ws.Cells.Item(530, 13).Value = 5000
$ws.Cells.Item(530, 14).Value = "`$/docena de matas"
$ws.Cells.Item(530, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(530, 16).Value = 833
$ws.Cells.Item(530, 17).Value = 6
$ws.Cells.Item(530, 18).Value = "Hortaliza"

$wb.Save()
